$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (Price column D, Volume(1h) column E).
# Column D values are prefixed with a literal leading apostrophe so Excel
# stores them as text (matching the source inlineStr cells) instead of
# re-parsing number-looking strings such as "161.80" into 161.8.

$ws.Range('D2').Value = '''30.262.66'
$ws.Range('E2').Value = '  +1.73%  '
$ws.Range('D3').Value = '''1.892.55'
$ws.Range('E3').Value = '  -1.24%  '
$ws.Range('D4').Value = '''1.002'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '''324.02'
$ws.Range('E5').Value = '  +1.98%  '
$ws.Range('E6').Value = '  +0.00%  '
$ws.Range('D7').Value = '''0.5184'
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  +1.19%  '
$ws.Range('D9').Value = '''0.08409'
$ws.Range('E9').Value = '  -1.34%  '
$ws.Range('E10').Value = '  +0.17%  '
$ws.Range('D11').Value = '''1.114'
$ws.Range('E11').Value = '  -0.82%  '
$ws.Range('D12').Value = '''23.19'
$ws.Range('E12').Value = '  +10.71%  '
$ws.Range('D13').Value = '''6.433'
$ws.Range('E13').Value = '  +2.07%  '
$ws.Range('D14').Value = '''1.892.98'
$ws.Range('E14').Value = '  -0.89%  '
$ws.Range('E15').Value = '  -0.65%  '
$ws.Range('E16').Value = '  +0.05%  '
$ws.Range('D17').Value = '''94.27'
$ws.Range('E17').Value = '  +0.17%  '
$ws.Range('D18').Value = '''0.00001108'
$ws.Range('E18').Value = '  -0.70%  '
$ws.Range('D19').Value = '''0.06639'
$ws.Range('E19').Value = '  -1.55%  '
$ws.Range('D20').Value = '''18.22'
$ws.Range('E20').Value = '  +1.40%  '
$ws.Range('E21').Value = '  +0.05%  '
$ws.Range('D22').Value = '''5.947'
$ws.Range('E22').Value = '  -1.53%  '
$ws.Range('D23').Value = '''30.239.84'
$ws.Range('E23').Value = '  +1.68%  '
$ws.Range('D24').Value = '''11.29'
$ws.Range('E24').Value = '  +0.64%  '
$ws.Range('E25').Value = '  +0.84%  '
$ws.Range('D26').Value = '''2.116.07'
$ws.Range('E26').Value = '  -0.57%  '
$ws.Range('D27').Value = '''21.62'
$ws.Range('E27').Value = '  +2.68%  '
$ws.Range('D28').Value = '''161.80'
$ws.Range('E28').Value = '  +1.61%  '
$ws.Range('D29').Value = '''2.330'
$ws.Range('E29').Value = '  -5.66%  '
$ws.Range('D30').Value = '''129.19'
$ws.Range('E30').Value = '  +0.27%  '
$ws.Range('D31').Value = '''1.087'
$ws.Range('E31').Value = '  +0.40%  '
$ws.Range('E32').Value = '  -0.45%  '
$ws.Range('D33').Value = '''6.094'
$ws.Range('E33').Value = '  -1.70%  '
$ws.Range('D34').Value = '''3.740'
$ws.Range('E34').Value = '  +1.63%  '
$ws.Range('D35').Value = '''0.02488'
$ws.Range('E35').Value = '  -0.60%  '
$ws.Range('D36').Value = '''0.06550'
$ws.Range('E36').Value = '  -1.45%  '
$ws.Range('D37').Value = '''5.326'
$ws.Range('E37').Value = '  +1.91%  '
$ws.Range('E38').Value = '  -0.42%  '
$ws.Range('E39').Value = '  -2.62%  '
$ws.Range('D40').Value = '''8.807'
$ws.Range('E40').Value = '  -3.53%  '
$ws.Range('D41').Value = '''11.82'
$ws.Range('E41').Value = '  +3.88%  '
$ws.Range('D42').Value = '''0.6494'
$ws.Range('E42').Value = '  -0.88%  '
$ws.Range('D43').Value = '''1.231'
$ws.Range('E43').Value = '  -0.93%  '
$ws.Range('D44').Value = '''0.6086'
$ws.Range('E44').Value = '  -0.67%  '
$ws.Range('E45').Value = '  +0.25%  '
$ws.Range('D46').Value = '''3.688'
$ws.Range('E46').Value = '  -0.03%  '
$ws.Range('D47').Value = '''2.054'
$ws.Range('E47').Value = '  -0.62%  '
$ws.Range('D48').Value = '''1.237'
$ws.Range('E48').Value = '  -0.34%  '
$ws.Range('D49').Value = '''124.76'
$ws.Range('E49').Value = '  -0.07%  '
$ws.Range('D50').Value = '''1.157'
$ws.Range('E50').Value = '  -2.61%  '
$ws.Range('D51').Value = '''79.07'
$ws.Range('E51').Value = '  +0.70%  '
